$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (A1:N1) ---
# A1 is left untouched (empty string cell, unchanged by this edit)
$ws.Range("B1").Value = "`$ bold('All')"
$ws.Range("C1").Value = "`$ bold('Europe')"
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Russia"
$ws.Range("M1").Value = "Saudi Arabia"
$ws.Range("N1").Value = "USA"

# --- Row labels (A2:A11) ---
$ws.Range("A2").Value = "Minimum tax of 2% on billionaires'`nwealth, in voluntary countries"
$ws.Range("A3").Value = "Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates"
$ws.Range("A4").Value = "L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage"
$ws.Range("A5").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("A6").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("A7").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("A8").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("A9").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("A10").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("A11").Value = "International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population"

# --- Numeric data (B2:N11) ---
$ws.Range("B2").Value = 0.746758615571303
$ws.Range("C2").Value = 0.899067107707109
$ws.Range("D2").Value = 0.88061965931131
$ws.Range("E2").Value = 0.859445533229057
$ws.Range("F2").Value = 0.964717436709091
$ws.Range("G2").Value = 0.667746196263347
$ws.Range("H2").Value = 0.853414072269436
$ws.Range("I2").Value = 0.937961080466205
$ws.Range("J2").Value = 0.775560245664184
$ws.Range("K2").Value = 0.52819369889068
$ws.Range("L2").Value = 0.661746313849443
$ws.Range("M2").Value = 0.844702008923008
$ws.Range("N2").Value = 0.619717771698186
$ws.Range("B3").Value = 0.53880897005617
$ws.Range("C3").Value = 0.622717294167845
$ws.Range("D3").Value = 0.472415613571565
$ws.Range("E3").Value = 0.553620391101682
$ws.Range("F3").Value = 0.834884450176832
$ws.Range("G3").Value = 0.28877771638812
$ws.Range("H3").Value = 0.625246524605311
$ws.Range("I3").Value = 0.740378276151058
$ws.Range("J3").Value = 0.54869733551346
$ws.Range("K3").Value = 0.370291814891383
$ws.Range("L3").Value = 0.622904912769984
$ws.Range("M3").Value = 0.834571283823028
$ws.Range("N3").Value = 0.39502177849613
$ws.Range("B4").Value = 0.457646793262474
$ws.Range("C4").Value = 0.487373574197681
$ws.Range("D4").Value = 0.339312817983837
$ws.Range("E4").Value = 0.399644861775946
$ws.Range("F4").Value = 0.61845489966271
$ws.Range("G4").Value = 0.36431897920809
$ws.Range("H4").Value = 0.63725872276991
$ws.Range("I4").Value = 0.470071392452268
$ws.Range("J4").Value = 0.336785160553625
$ws.Range("K4").Value = 0.260444114788849
$ws.Range("L4").Value = 0.695125891930721
$ws.Range("M4").Value = 0.958227224438065
$ws.Range("N4").Value = 0.312014305831053
$ws.Range("B5").Value = 0.339613614975512
$ws.Range("C5").Value = 0.340947849011829
$ws.Range("D5").Value = 0.23045269222624
$ws.Range("E5").Value = 0.244638275828069
$ws.Range("F5").Value = 0.466259178198866
$ws.Range("G5").Value = 0.0995292164645176
$ws.Range("H5").Value = 0.557119912781717
$ws.Range("I5").Value = 0.365610055382771
$ws.Range("J5").Value = 0.290999563900097
$ws.Range("K5").Value = 0.0719327097269795
$ws.Range("L5").Value = 0.658406804039628
$ws.Range("M5").Value = 0.808599687276629
$ws.Range("N5").Value = 0.200736276966102
$ws.Range("B6").Value = 0.338691995923963
$ws.Range("C6").Value = 0.366925945043564
$ws.Range("D6").Value = 0.207521570407308
$ws.Range("E6").Value = 0.104264439406356
$ws.Range("F6").Value = 0.604095250398998
$ws.Range("G6").Value = 0.424311899764081
$ws.Range("H6").Value = 0.474036769409407
$ws.Range("I6").Value = 0.473707266132616
$ws.Range("J6").Value = 0.269478210110152
$ws.Range("K6").Value = 0.18547588627103
$ws.Range("L6").Value = 0.409824606520925
$ws.Range("M6").Value = 0.932399297563257
$ws.Range("N6").Value = 0.263563439563823
$ws.Range("B7").Value = 0.329110709214005
$ws.Range("C7").Value = 0.523687810073322
$ws.Range("D7").Value = 0.42064295930788
$ws.Range("E7").Value = 0.446717126447487
$ws.Range("F7").Value = 0.844410613388007
$ws.Range("G7").Value = 0.2309618225646
$ws.Range("H7").Value = 0.428000591337885
$ws.Range("I7").Value = 0.491629330123401
$ws.Range("J7").Value = 0.266420209147519
$ws.Range("K7").Value = 0.328888261683248
$ws.Range("L7").Value = -0.0778417618336274
$ws.Range("M7").Value = 0.497163901847429
$ws.Range("N7").Value = 0.201989188661513
$ws.Range("B8").Value = 0.320103458109216
$ws.Range("C8").Value = 0.385513480941771
$ws.Range("D8").Value = 0.263621500427078
$ws.Range("E8").Value = 0.362672359948726
$ws.Range("F8").Value = 0.466268824250188
$ws.Range("G8").Value = 0.187158200316393
$ws.Range("H8").Value = 0.528868523638982
$ws.Range("I8").Value = 0.382752450036814
$ws.Range("J8").Value = 0.28229502195541
$ws.Range("K8").Value = 0.044077194018155
$ws.Range("L8").Value = 0.697267347192898
$ws.Range("M8").Value = 0.776081506237869
$ws.Range("N8").Value = 0.109587709651655
$ws.Range("B9").Value = 0.31788661467437
$ws.Range("C9").Value = 0.420490345740477
$ws.Range("D9").Value = 0.49292005022254
$ws.Range("E9").Value = 0.304386935822869
$ws.Range("F9").Value = 0.499615146146505
$ws.Range("G9").Value = 0.0892781446390868
$ws.Range("H9").Value = 0.456488900986814
$ws.Range("I9").Value = 0.40971182484094
$ws.Range("J9").Value = 0.38981506699751
$ws.Range("K9").Value = 0.0340824452221103
$ws.Range("L9").Value = 0.36590464407975
$ws.Range("M9").Value = 0.660342805051329
$ws.Range("N9").Value = 0.20372593232651
$ws.Range("B10").Value = 0.315567059102638
$ws.Range("C10").Value = 0.492563434872145
$ws.Range("D10").Value = 0.347566986960034
$ws.Range("E10").Value = 0.455512915855029
$ws.Range("F10").Value = 0.569963435571414
$ws.Range("G10").Value = 0.285573393035941
$ws.Range("H10").Value = 0.505538471803289
$ws.Range("I10").Value = 0.493970781786548
$ws.Range("J10").Value = 0.372640046923682
$ws.Range("K10").Value = 0.21134660742927
$ws.Range("L10").Value = -0.0291282837527572
$ws.Range("M10").Value = 0.707615721886186
$ws.Range("N10").Value = 0.209230862824192
$ws.Range("B11").Value = 0.00973878914166017
$ws.Range("C11").Value = 0.0839217421719147
$ws.Range("D11").Value = 0.166213268172324
$ws.Range("E11").Value = 0.0835668003633305
$ws.Range("F11").Value = -0.00838446613759675
$ws.Range("G11").Value = -0.0810140794603717
$ws.Range("H11").Value = -0.00387824155049399
$ws.Range("I11").Value = 0.0370952506022134
$ws.Range("J11").Value = 0.00387057722535984
$ws.Range("K11").Value = -0.105598258636374
$ws.Range("L11").Value = -0.0180680286245144
$ws.Range("M11").Value = 0.414946008655522
$ws.Range("N11").Value = -0.100891956929632
